# Fix Juniper infrastructure costs - remove Professional Services category
# - Remove all Professional Services line items from Infrastructure Costs,
#   Credits, and 3-Year Summary sheets.
# - Rename "Software" category label to "Software Licenses" everywhere it
#   is used as a category tag.
# - Bump the Cover sheet "Solution" (Generated) date.
# - Keep the _xlnm._FilterDatabase defined names and each sheet's AutoFilter
#   range in sync with the new (smaller) data extents.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Cover sheet: bump the date in C7
# ---------------------------------------------------------------------
$wsCover = $wb.Worksheets.Item("Cover")
$wsCover.Range("C7").Value = "November 24, 2025"

# ---------------------------------------------------------------------
# 2. Infrastructure Costs sheet
# ---------------------------------------------------------------------
$wsInfra = $wb.Worksheets.Item("Infrastructure Costs")

# Rename "Software" -> "Software Licenses" for rows 7-12 (category column A)
foreach ($r in 7..12) {
    $wsInfra.Cells.Item($r, 1).Value = "Software Licenses"
}

# Remove the 9 "Professional Services" line items (rows 15-23). The TOTAL
# row (formerly row 24) shifts up to row 15 and its SUM() formulas are
# re-based automatically by the row delete.
$wsInfra.Range("A15:A23").EntireRow.Delete()

# Resync AutoFilter range + _xlnm._FilterDatabase defined name to A2:K15
$wsInfra.AutoFilterMode = $false
[void]$wsInfra.Range("A2:K15").AutoFilter()
foreach ($n in $wb.Names) {
    if ($n.Name -eq "Infrastructure Costs!_FilterDatabase") {
        $n.RefersTo = "='Infrastructure Costs'!`$A`$2:`$K`$15"
    }
}

# ---------------------------------------------------------------------
# 3. Credits sheet
# ---------------------------------------------------------------------
$wsCredits = $wb.Worksheets.Item("Credits")

# Rename "Software" -> "Software Licenses" for row 4 (category column A)
$wsCredits.Cells.Item(4, 1).Value = "Software Licenses"

# Remove the "Professional Services" / Partner Credit row (row 6)
$wsCredits.Rows.Item(6).Delete()

# Resync AutoFilter range + _xlnm._FilterDatabase defined name to A2:D5
$wsCredits.AutoFilterMode = $false
$wsCredits.Range("A2:D5").AutoFilter()
foreach ($n in $wb.Names) {
    if ($n.Name -eq "Credits!_FilterDatabase") {
        $n.RefersTo = "=Credits!`$A`$2:`$D`$5"
    }
}

# ---------------------------------------------------------------------
# 4. 3-Year Summary sheet
# ---------------------------------------------------------------------
$wsSummary = $wb.Worksheets.Item("3-Year Summary")

# Rename "Software" -> "Software Licenses" for row 4 (category column A)
$wsSummary.Cells.Item(4, 1).Value = "Software Licenses"

# Remove the "Professional Services" row (row 6). The TOTAL row (formerly
# row 7) shifts up to row 6 and its SUM() formulas are re-based
# automatically by the row delete.
$wsSummary.Rows.Item(6).Delete()

# Resync AutoFilter range + _xlnm._FilterDatabase defined name to A2:G6
$wsSummary.AutoFilterMode = $false
$wsSummary.Range("A2:G6").AutoFilter()
foreach ($n in $wb.Names) {
    if ($n.Name -eq "3-Year Summary!_FilterDatabase") {
        $n.RefersTo = "='3-Year Summary'!`$A`$2:`$G`$6"
    }
}
